$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the full data row (columns B..AC) between paired rows. ---
# (Column A, the running index, stays put; everything else trades places.)
$rowPairs = @(
    @(21,22),
    @(24,25),
    @(37,38),
    @(42,43),
    @(45,46),
    @(50,51),
    @(54,55),
    @(60,61),
    @(75,76),
    @(78,79)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("B" + $r1 + ":AC" + $r1)
    $rng2 = $ws.Range("B" + $r2 + ":AC" + $r2)

    for ($c = 2; $c -le 29; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# --- Append the newly-scraped fixtures as rows 80-85. ---
function Set-Row {
    param($r, $a, $b, $f, $g, $e, $k, $l, $m, $n, $o, $p, $q, $rr, $s, $t, $u, $v, $w, $x, $y, $z, $aa)

    $ws.Cells.Item($r, 1).Value2  = $a
    $ws.Cells.Item($r, 2).Value2  = $b
    $ws.Cells.Item($r, 3).Value2  = "Qatar Stars League"
    $ws.Cells.Item($r, 4).Value2  = "Qatar Stars League"
    $ws.Cells.Item($r, 5).Value2  = $e
    $ws.Cells.Item($r, 6).Value2  = $f
    $ws.Cells.Item($r, 7).Value2  = $g
    $ws.Cells.Item($r, 11).Value2 = $k
    $ws.Cells.Item($r, 12).Value2 = $l
    $ws.Cells.Item($r, 13).Value2 = $m
    $ws.Cells.Item($r, 14).Value2 = $n
    $ws.Cells.Item($r, 15).Value2 = $o
    $ws.Cells.Item($r, 16).Value2 = $p
    $ws.Cells.Item($r, 17).Value2 = $q
    $ws.Cells.Item($r, 18).Value2 = $rr
    $ws.Cells.Item($r, 19).Value2 = $s
    $ws.Cells.Item($r, 20).Value2 = $t
    $ws.Cells.Item($r, 21).Value2 = $u
    $ws.Cells.Item($r, 22).Value2 = $v
    $ws.Cells.Item($r, 23).Value2 = $w
    $ws.Cells.Item($r, 24).Value2 = $x
    $ws.Cells.Item($r, 25).Value2 = $y
    $ws.Cells.Item($r, 26).Value2 = $z
    $ws.Cells.Item($r, 27).Value2 = $aa

    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160
    $ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Set-Row 80 78 7840801 "AlShamal SC"  "AlMuaidar"   45351.5           2.25 4    2.4  2.4   4    2.25  0     1.95  1.85 3    1.9   1.9   0 0 0 0 0
Set-Row 81 79 7840805 "AlWakrah SC" "Al Markhiya" 45351.58333333334 1.062 11  17   1.25  5.5  10    -1.5  1.825 1.975 3.25 1.875 1.925 0 0 0 0 0
Set-Row 82 80 7840802 "AlAhli Doha" "Umm Salal"   45351.58333333334 2.4   4   2.25 2.2   4    2.5   0     1.825 1.975 3    1.825 1.975 0 0 0 0 0
Set-Row 83 81 7840803 "Al Sadd"     "Qatar SC Doha" 45352.5         1.285 5.75 7.5 1.285 5.75 7.5   -1.75 1.975 1.825 3.5  1.95  1.85  0 0 0 0 0
Set-Row 84 82 7840685 "Al Duhail"   "Al Gharafa"  45352.5           2.1   3.8 2.9  1.909 4    3.3   -0.5  1.925 1.875 3.5  1.9   1.9   0 0 0 0 0
Set-Row 85 83 7840804 "AlArabi Doha" "AlRayyan SC" 45352.58333333334 2.7  3.75 2.2  2.5   3.75 2.375 0     1.975 1.825 3.25 1.95  1.85  0 0 0 0 0
